$wb = $excel.ActiveWorkbook

# Rename Sheet1 to DisplayValues
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "DisplayValues"

# Add the new named ranges (workbook scope)
$wb.Names.Add("disp_value_col_head", "='DisplayValues'!`$B`$1:`$C`$1")
$wb.Names.Add("disp_value_row_head", "='DisplayValues'!`$A`$2:`$A`$3")
$wb.Names.Add("dis_value_values", "='DisplayValues'!`$B`$2:`$C`$3")
$wb.Names.Add("cuts", "='Lookups'!`$A`$1:`$A`$2")
$wb.Names.Add("cuts_config", "='Lookups'!`$A`$1:`$E`$2")
$wb.Names.Add("default_menu", "='Lookups'!`$E`$2:`$E`$101")
$wb.Names.Add("cuts_head", "='Lookups'!`$F`$1:`$G`$1")

$wb.Save()
